$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.217.76"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "3.665.06"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'595.85"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "'192.23"
$ws.Range("E6").Value = "  +6.69%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.698"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("D10").Value = "'0.153"
$ws.Range("E10").Value = "  -5.26%  "
$ws.Range("D11").Value = "'57.13"
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("D12").Value = "'0.0000272"
$ws.Range("E12").Value = "  -5.59%  "
$ws.Range("D13").Value = "'10.24"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "4.254.73"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").Value = "3.671.91"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "'18.87"
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "68.047.63"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'12.60"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").Value = "'404.72"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "'4.41"
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("D23").Value = "'88.34"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "'2.95"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'12.59"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "'10.87"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").Value = "'6.07"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'3.74"
$ws.Range("E28").Value = "  -9.46%  "
$ws.Range("D29").Value = "'9.35"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").Value = "'31.97"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("D31").Value = "'7.17"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").Value = "'12.28"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("D33").Value = "'67.14"
$ws.Range("E33").Value = "  +4.90%  "
$ws.Range("D34").Value = "'43.94"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("D35").Value = "'0.116"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'607.30"
$ws.Range("E36").Value = "  +3.00%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'0.392"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "0.0₃0774"
$ws.Range("E40").Value = "  -11.13%  "
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").Value = "'2.90"
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("D43").Value = "'0.0426"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").Value = "'2.51"
$ws.Range("E44").Value = "  -9.05%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.22"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.136"
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("D47").Value = "2.777.98"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").Value = "'8.93"
$ws.Range("E48").Value = "  -2.71%  "
$ws.Range("D49").Value = "'143.90"
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("E50").Value = "  -4.35%  "
$ws.Range("D51").Value = "'2.54"
$ws.Range("E51").Value = "  -12.01%  "
